# Adds 4 new test cases to the "Admin" sheet (rows 9-16), mirroring the
# existing TC04_Admin_CreateJobTitle block (rows 7-8): a "header" row with
# the test-case name + field names, followed by a "data" row with the
# corresponding values and a trailing back-reference to the test-case name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")

# ---------------------------------------------------------------------
# Helper: copy the *format only* of a reference cell onto a target cell,
# so we reuse existing cellXfs entries instead of inventing new ones.
# ---------------------------------------------------------------------
function Copy-Format($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$hdrA = $ws.Range("A7")   # style 12 - header test-case name (col A)
$hdrB = $ws.Range("B7")   # style 11 - header field name
$refC = $ws.Range("D8")   # style 12 - data row back-reference to TC name
$lblData = $ws.Range("B8")# style 10 - data row value cell
$emptyData = $ws.Range("D7") # style 10 - empty styled cell (data row, col A/C)

# ===================== TC05_Admin_CreateEmploymentStatus =====================
Copy-Format $hdrA $ws.Range("A9")
$ws.Range("A9").Value2 = "TC05_Admin_CreateEmploymentStatus"

Copy-Format $hdrB $ws.Range("B9")
$ws.Range("B9").Value2 = "EmploymentStatus"

Copy-Format $emptyData $ws.Range("C9")

Copy-Format $ws.Range("C2") $ws.Range("D9")   # style 1 (no value)

Copy-Format $emptyData $ws.Range("A10")

Copy-Format $lblData $ws.Range("B10")
$ws.Range("B10").Value2 = "Full-Time Internship"

Copy-Format $hdrA $ws.Range("C10")
$ws.Range("C10").Value2 = "TC05_Admin_CreateEmploymentStatus"

# ===================== TC06_Admin_CreateJobCategory =====================
Copy-Format $hdrA $ws.Range("A11")
$ws.Range("A11").Value2 = "TC06_Admin_CreateJobCategory"

Copy-Format $hdrB $ws.Range("B11")
$ws.Range("B11").Value2 = "JobCategory"

Copy-Format $emptyData $ws.Range("C11")

Copy-Format $emptyData $ws.Range("A12")

Copy-Format $lblData $ws.Range("B12")
$ws.Range("B12").Value2 = "Search And Online Marketing"

Copy-Format $hdrA $ws.Range("C12")
$ws.Range("C12").Value2 = "TC06_Admin_CreateJobCategory"

# ===================== TC07_Admin_CreateWorkShift =====================
Copy-Format $hdrA $ws.Range("A13")
$ws.Range("A13").Value2 = "TC07_Admin_CreateWorkShift"

Copy-Format $hdrB $ws.Range("B13")
$ws.Range("B13").Value2 = "WorkShiftName"

Copy-Format $hdrB $ws.Range("C13")
$ws.Range("C13").Value2 = "FromTime"

Copy-Format $hdrB $ws.Range("D13")
$ws.Range("D13").Value2 = "ToTime"

Copy-Format $emptyData $ws.Range("A14")

Copy-Format $lblData $ws.Range("B14")
$ws.Range("B14").Value2 = "Morning Shift"

# C14: text "08:00 AM" formatted with the time number format + quote prefix
$ws.Range("C14").Value2 = "'08:00 AM"
$ws.Range("C14").NumberFormat = "h:mm AM/PM"

# D14: text "05:00 PM" stored with a quote-prefix (General format)
$ws.Range("D14").Value2 = "'05:00 PM"

Copy-Format $hdrA $ws.Range("E14")
$ws.Range("E14").Value2 = "TC07_Admin_CreateWorkShift"

# ===================== TC08_Admin_CreatePayGrade =====================
Copy-Format $hdrA $ws.Range("A15")
$ws.Range("A15").Value2 = "TC08_Admin_CreatePayGrade"

Copy-Format $hdrB $ws.Range("B15")
$ws.Range("B15").Value2 = "PayGrade"

Copy-Format $hdrB $ws.Range("C15")
$ws.Range("C15").Value2 = "Currency"

Copy-Format $hdrB $ws.Range("D15")
$ws.Range("D15").Value2 = "MinSalary"

Copy-Format $hdrB $ws.Range("E15")
$ws.Range("E15").Value2 = "MaxSalary"

Copy-Format $emptyData $ws.Range("A16")

Copy-Format $lblData $ws.Range("B16")
$ws.Range("B16").Value2 = "Grade 10"

# C16: text "EUR - Euro" but reuses the same quote-prefixed time-format style
$ws.Range("C16").Value2 = "'EUR - Euro"
$ws.Range("C16").NumberFormat = "h:mm AM/PM"

# D16 / E16: text "10" / "1000" stored with a quote-prefix (General format)
$ws.Range("D16").Value2 = "'10"
$ws.Range("E16").Value2 = "'1000"

Copy-Format $hdrA $ws.Range("F16")
$ws.Range("F16").Value2 = "TC08_Admin_CreatePayGrade"

# ---------------------------------------------------------------------
# View state: scroll the sheet down a bit and move the lingering
# out-of-range selection from J18 to C18.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
